$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/51f6ca464ded53573595c3f7f297d196a35b516c/e2e/7fc1593d-971e-4343-8a7a-3c0207b7150e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec0ad3746077f3590ddb3fd51a23d8645553f0ee/e2e/7fc1593d-971e-4343-8a7a-3c0207b7150e.md."

# Overview sheet - row 3 is the 7fc1593d file
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-27 02:47:45"

# zh-cn sheet - row 3
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-27 02:47:40"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet - row 3
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-27 02:47:45"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
